$wb = $excel.ActiveWorkbook

# Sheets that contain the full event listing table: "展览" (sheet1) and "全部类型" (sheet4)
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 (the oldest/cancelled event) is removed; rows 3-37 shift up to rows 2-36.
    $ws.Rows.Item(2).Delete()

    # Column A is a sequential index (0-based) that must stay 1..35 for rows 2..36
    # after the shift (Excel's row delete shifts A along with everything else, so reset it).
    for ($i = 2; $i -le 36; $i++) {
        $ws.Cells.Item($i, 1).Value = $i - 1
    }

    # "想去人数" (F) / "最低票价" (G) are live counters re-scraped at commit time; they
    # do not line up with a pure shift of the row below, so set them explicitly.
    $ws.Range("F2").Value = 190
    $ws.Range("G2").Value = 50
    $ws.Range("F3").Value = 2973
    $ws.Range("G3").Value = 65
    $ws.Range("F4").Value = 211
    $ws.Range("G4").Value = 55
    $ws.Range("F5").Value = 111
    $ws.Range("G5").Value = 55
    $ws.Range("F6").Value = 191
    $ws.Range("G6").Value = 48
    $ws.Range("F7").Value = 1636
    $ws.Range("G7").Value = 60
    $ws.Range("F8").Value = 1605
    $ws.Range("G8").Value = 55
    $ws.Range("F9").Value = 53
    $ws.Range("G9").Value = 45
    $ws.Range("F10").Value = 351
    $ws.Range("G10").Value = 55
    $ws.Range("F11").Value = 234
    $ws.Range("G11").Value = "不可售"
    $ws.Range("F12").Value = 27
    $ws.Range("G12").Value = 22.33
    $ws.Range("F13").Value = 186
    $ws.Range("G13").Value = 55
    $ws.Range("F14").Value = 25
    $ws.Range("G14").Value = 45
    $ws.Range("F15").Value = 223
    $ws.Range("G15").Value = 55
    $ws.Range("F16").Value = 231
    $ws.Range("G16").Value = 55
    $ws.Range("F17").Value = 226
    $ws.Range("G17").Value = 55
    $ws.Range("F18").Value = 20
    $ws.Range("G18").Value = 55
    $ws.Range("F19").Value = 11
    $ws.Range("G19").Value = 55
    $ws.Range("F20").Value = 37
    $ws.Range("G20").Value = 55
    $ws.Range("F21").Value = 5
    $ws.Range("G21").Value = 55
    $ws.Range("F22").Value = 350
    $ws.Range("G22").Value = 52.1
    $ws.Range("F23").Value = 141
    $ws.Range("G23").Value = 55
    $ws.Range("F24").Value = 93
    $ws.Range("G24").Value = 52.5
    $ws.Range("F25").Value = 16
    $ws.Range("G25").Value = 40
    $ws.Range("F26").Value = 1986
    $ws.Range("G26").Value = 69
    $ws.Range("F27").Value = 50
    $ws.Range("G27").Value = 56
    $ws.Range("F28").Value = 453
    $ws.Range("G28").Value = 64
    $ws.Range("F29").Value = 12
    $ws.Range("G29").Value = 45
    $ws.Range("F30").Value = 165
    $ws.Range("G30").Value = 55
    $ws.Range("F31").Value = 570
    $ws.Range("G31").Value = "已售罄"
    $ws.Range("F32").Value = 221
    $ws.Range("G32").Value = 45
    $ws.Range("F33").Value = 328
    $ws.Range("G33").Value = 55
    $ws.Range("F34").Value = 3
    $ws.Range("G34").Value = 45
    $ws.Range("F35").Value = 486
    $ws.Range("G35").Value = 45
    $ws.Range("F36").Value = 4
    $ws.Range("G36").Value = 45
}

# F24 ("南昌·漫拥动漫嘉年华Pro-追光启航" want-to-go count) differs very slightly between
# the two sheets in the source snapshot used for this commit.
$wb.Worksheets.Item("展览").Range("F24").Value = 93
$wb.Worksheets.Item("全部类型").Range("F24").Value = 94

Write-Output "done"
